$wb = $excel.ActiveWorkbook

# "sign up" sheet (index 2) - faker data updates
$ws = $wb.Worksheets.Item(2)

# Row 2 (Bloom Bee): phone number changed, email text changed
$ws.Range("B2").Value = 9800321453
$ws.Range("C2").Value = "bee1@gmail.com"

# Row 3: phone number replaced with a large faker-generated numeric value
$ws.Range("B3").Value = 9825088978687871000.0

# Row 9 (Margaret Ramirez DVM): phone number becomes a text value
$ws.Range("B9").Value = "9772..552113"

# Make "sign up" the active sheet/tab, with cell C8 selected
$ws.Activate() | Out-Null
$ws.Range("C8").Select() | Out-Null
